$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,3,3,0),
    @(3,1,4,2),
    @(6,0,5,3),
    @(4,1,4,2),
    @(5,0,5,2),
    @(7,3,6,0),
    @(2,2,3,1),
    @(5,0,3,3),
    @(4,0,4,3),
    @(4,0,6,3),
    @(4,0,3,3),
    @(4,2,4,0),
    @(3,0,3,3),
    @(3,1,3,2),
    @(6,2,4,1),
    @(3,3,3,0),
    @(6,2,7,0),
    @(4,0,5,3),
    @(3,0,3,3),
    @(3,0,3,3),
    @(6,1,7,2),
    @(5,2,5,0),
    @(6,0,4,2),
    @(4,3,5,0),
    @(4,2,4,1),
    @(5,0,7,3),
    @(3,1,3,2),
    @(3,2,3,1),
    @(7,2,6,0),
    @(5,3,4,0),
    @(4,1,5,2),
    @(4,2,4,1),
    @(5,0,5,2),
    @(6,2,5,0),
    @(4,2,6,0),
    @(4,3,3,0)
)

$startRow = 2556
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$ws.Range("A2592").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2572
